$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # "CashFlowPlan"
$ws2 = $wb.Worksheets.Item(2)   # "Definition"

# -------------------------------------------------------------------
# Sheet1 "CashFlowPlan": insert a new column C ("BaseAmount") and
# rename the old "Amount" header to "LocalAmount". The previous
# column C ("Notes") shifts to column D.
# -------------------------------------------------------------------
$ws1.Range("C1").EntireColumn.Insert() | Out-Null

$ws1.Range("B1").Value = "LocalAmount"
$ws1.Range("C1").Value = "BaseAmount"
$ws1.Columns("A:D").AutoFit() | Out-Null

# -------------------------------------------------------------------
# Sheet2 "Definition": split the "Amount" definition row into two
# rows - "LocalAmount" and "BaseAmount" - and push the "Notes" row
# (and the legend rows below it) down by one row.
# -------------------------------------------------------------------
$ws2.Range("A9").EntireRow.Insert() | Out-Null

# Row 8 now describes "LocalAmount" (was "Amount"/"M")
$ws2.Range("B8").Value = "LocalAmount"
$ws2.Range("D8").Value = "Local Amount"
$ws2.Range("E8").Value = "O"

# Row 9 (newly inserted) describes "BaseAmount"
$ws2.Range("A9").Value = 3
$ws2.Range("B9").Value = "BaseAmount"
$ws2.Range("D9").Value = "Base Amount"
$ws2.Range("E9").Value = "O"
$ws2.Range("F9").Value = "Numeric"
$ws2.Range("G9").Value = "18,2"

# -------------------------------------------------------------------
# Restore the selections recorded in each sheet view.
# -------------------------------------------------------------------
$ws2.Activate() | Out-Null
$ws2.Range("E9").Select() | Out-Null

$ws1.Activate() | Out-Null
$ws1.Range("B5").Select() | Out-Null
